$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1362.3
$ws.Range("J29").Value = 2219.8333
$ws.Range("L29").Value = 6659.499899999999
$ws.Range("N29").Value = -7221.499899999999
$ws.Range("H32").Value = 21249.25
$ws.Range("I32").Value = 26665.666
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 26665.666
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -26339.666
$ws.Range("N32").Value = -5652
$ws.Range("H38").Value = 787.375
$ws.Range("I38").Value = 666.3333
$ws.Range("J38").Value = 860
$ws.Range("K38").Value = 1998.9999
$ws.Range("L38").Value = 2580
$ws.Range("M38").Value = -1626.9999
$ws.Range("N38").Value = -3324
$ws.Range("H53").Value = 132.52632
$ws.Range("I53").Value = 114.92857
$ws.Range("J53").Value = 181.8
$ws.Range("K53").Value = 114.92857
$ws.Range("L53").Value = 181.8
$ws.Range("M53").Value = 522.07143
$ws.Range("N53").Value = -1455.8
$ws.Range("H98").Value = 1312.125
$ws.Range("I98").Value = 1271.1428
$ws.Range("K98").Value = 1271.1428
$ws.Range("M98").Value = 226.8571999999999
$ws.Range("H116").Value = 3239.6
$ws.Range("I116").Value = 2549.5
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 2549.5
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 892.5
$ws.Range("N116").Value = -12884
$ws.Range("H122").Value = 1312.125
$ws.Range("I122").Value = 1271.1428
$ws.Range("K122").Value = 3813.4284
$ws.Range("M122").Value = -1363.4284
$ws.Range("H132").Value = 2321.5
$ws.Range("I132").Value = 1260.2727
$ws.Range("K132").Value = 3780.8181
$ws.Range("M132").Value = -1250.8181
$ws.Range("H138").Value = 2018.6349
$ws.Range("I138").Value = 1638.25
$ws.Range("J138").Value = 2148.1277
$ws.Range("K138").Value = 4914.75
$ws.Range("L138").Value = 6444.3831
$ws.Range("M138").Value = 225.25
$ws.Range("N138").Value = -16724.3831
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 20214.285
$ws.Range("I37").Value = 7500
$ws.Range("J37").Value = 22333.334
$ws.Range("K37").Value = 7500
$ws.Range("L37").Value = 22333.334
$ws.Range("M37").Value = -7227
$ws.Range("N37").Value = -22879.334
$ws.Range("H45").Value = 2416.0715
$ws.Range("I45").Value = 1742.5
$ws.Range("K45").Value = 1742.5
$ws.Range("M45").Value = -1365.5
$ws.Range("H61").Value = 6242.125
$ws.Range("I61").Value = 6242.125
$ws.Range("K61").Value = 6242.125
$ws.Range("M61").Value = -6030.125
$ws.Range("H63").Value = 2100
$ws.Range("I63").Value = 1900
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 1900
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -1214
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2100
$ws.Range("I66").Value = 1900
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 9500
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -6068
$ws.Range("N66").Value = -19364
$ws.Range("H74").Value = 1422.8334
$ws.Range("J74").Value = 2402.5
$ws.Range("L74").Value = 2402.5
$ws.Range("N74").Value = -4150.5
$ws.Range("H77").Value = 1422.8334
$ws.Range("J77").Value = 2402.5
$ws.Range("L77").Value = 12012.5
$ws.Range("N77").Value = -20748.5
$ws.Range("H122").Value = 11808.036
$ws.Range("I122").Value = 8524.038
$ws.Range("K122").Value = 25572.114
$ws.Range("M122").Value = -23122.114
$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820
$ws.Range("H132").Value = 4082.8
$ws.Range("I132").Value = 3887.4119
$ws.Range("K132").Value = 11662.2357
$ws.Range("M132").Value = -9132.235700000001
$ws.Range("H136").Value = 6242.125
$ws.Range("I136").Value = 6242.125
$ws.Range("K136").Value = 18726.375
$ws.Range("M136").Value = -16176.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29999.5
$ws.Range("J35").Value = 29999.5
$ws.Range("L35").Value = 29999.5
$ws.Range("N35").Value = -30619.5
$ws.Range("H96").Value = 19999.334
$ws.Range("I96").Value = 19999.334
$ws.Range("K96").Value = 19999.334
$ws.Range("M96").Value = -17253.334
$ws.Range("H105").Value = 2735.4375
$ws.Range("I105").Value = 2066.8462
$ws.Range("K105").Value = 2066.8462
$ws.Range("M105").Value = -319.8462
$ws.Range("H134").Value = 3167.6667
$ws.Range("I134").Value = 2813.625
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 8440.875
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -5905.875
$ws.Range("N134").Value = -23070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 83812.19
$ws.Range("I22").Value = 132777.67
$ws.Range("J22").Value = 20856.572
$ws.Range("K22").Value = 132777.67
$ws.Range("L22").Value = 20856.572
$ws.Range("M22").Value = -132427.67
$ws.Range("N22").Value = -21556.572
$ws.Range("H41").Value = 17373.5
$ws.Range("I41").Value = 4500
$ws.Range("J41").Value = 21664.666
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 21664.666
$ws.Range("M41").Value = -4072
$ws.Range("N41").Value = -22520.666
$ws.Range("H105").Value = 881.0769
$ws.Range("I105").Value = 881.0769
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 881.0769
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 865.9231
$ws.Range("N105").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 250412.25
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 500274.5
$ws.Range("K107").Value = 1650
$ws.Range("L107").Value = 1500823.5
$ws.Range("M107").Value = 270
$ws.Range("N107").Value = -1504663.5
$ws.Range("H132").Value = 2867.8
$ws.Range("I132").Value = 2867.8
$ws.Range("K132").Value = 25810.2
$ws.Range("M132").Value = -23280.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 4
$ws.Range("I13").Value = 4
$ws.Range("K13").Value = 4
$ws.Range("M13").Value = 135
$ws.Range("H18").Value = 1846134.4
$ws.Range("J18").Value = 19200.25
$ws.Range("L18").Value = 19200.25
$ws.Range("N18").Value = -19786.25
$ws.Range("H21").Value = 25000006
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H30").Value = 25000006
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H80").Value = 2099.8
$ws.Range("I80").Value = 2499.75
$ws.Range("K80").Value = 2499.75
$ws.Range("M80").Value = -1501.75
$ws.Range("H83").Value = 2099.8
$ws.Range("I83").Value = 2499.75
$ws.Range("K83").Value = 12498.75
$ws.Range("M83").Value = -7506.75
$ws.Range("H122").Value = 69293.734
$ws.Range("I122").Value = 1948.9
$ws.Range("J122").Value = 203983.4
$ws.Range("K122").Value = 5846.700000000001
$ws.Range("L122").Value = 611950.2
$ws.Range("M122").Value = -3396.700000000001
$ws.Range("N122").Value = -616850.2
$ws.Range("H132").Value = 2418.375
$ws.Range("I132").Value = 2571.6572
$ws.Range("K132").Value = 7714.971600000001
$ws.Range("M132").Value = -5184.971600000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1250
$ws.Range("I7").Value = 1250
$ws.Range("K7").Value = 1250
$ws.Range("M7").Value = -1138
$ws.Range("H22").Value = 3293.0881
$ws.Range("I22").Value = 1845.0588
$ws.Range("J22").Value = 4741.1177
$ws.Range("K22").Value = 1845.0588
$ws.Range("L22").Value = 4741.1177
$ws.Range("M22").Value = -1550.0588
$ws.Range("N22").Value = -5331.1177
$ws.Range("H27").Value = 3293.0881
$ws.Range("I27").Value = 1845.0588
$ws.Range("J27").Value = 4741.1177
$ws.Range("K27").Value = 1845.0588
$ws.Range("L27").Value = 4741.1177
$ws.Range("M27").Value = -1738.0588
$ws.Range("N27").Value = -4955.1177
$ws.Range("H40").Value = 3020.4333
$ws.Range("I40").Value = 2969.4644
$ws.Range("K40").Value = 2969.4644
$ws.Range("M40").Value = -2833.4644
$ws.Range("H46").Value = 3789.4211
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 5499.8335
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 5499.8335
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -5875.8335
$ws.Range("H82").Value = 112216.555
$ws.Range("I82").Value = 1135.7142
$ws.Range("J82").Value = 500999.5
$ws.Range("K82").Value = 1135.7142
$ws.Range("L82").Value = 500999.5
$ws.Range("M82").Value = -774.7141999999999
$ws.Range("N82").Value = -501721.5
$ws.Range("H85").Value = 112216.555
$ws.Range("I85").Value = 1135.7142
$ws.Range("J85").Value = 500999.5
$ws.Range("K85").Value = 1135.7142
$ws.Range("L85").Value = 500999.5
$ws.Range("M85").Value = 112.2858000000001
$ws.Range("N85").Value = -503495.5
$ws.Range("H122").Value = 6036.7144
$ws.Range("J122").Value = 6994.933
$ws.Range("L122").Value = 20984.799
$ws.Range("N122").Value = -25884.799
$ws.Range("H126").Value = 1250
$ws.Range("I126").Value = 1250
$ws.Range("K126").Value = 3750
$ws.Range("M126").Value = -1280
$ws.Range("H132").Value = 2107.1428
$ws.Range("I132").Value = 1750
$ws.Range("K132").Value = 5250
$ws.Range("M132").Value = -2720
$ws.Range("H134").Value = 65000
$ws.Range("J134").Value = 65000
$ws.Range("L134").Value = 65000
$ws.Range("N134").Value = -75140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 45000
$ws.Range("J24").Value = 45000
$ws.Range("L24").Value = 45000
$ws.Range("N24").Value = -45460
$ws.Range("H26").Value = 516875
$ws.Range("J26").Value = 516875
$ws.Range("L26").Value = 516875
$ws.Range("N26").Value = -517461
